# Apply weekly Chirimoya price update: insert 3 new rows (new week, 2023-10-13 /
# serial 45212) right before the existing row 352 block, shifting the rest of
# the data (old rows 352:371) down to 355:374.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above current row 352; existing rows 352:371 shift down
# to 355:374 and inherit formatting from the row above the insertion point
# automatically with InsertBefore semantics, so we set formatting explicitly
# afterwards to be safe.
$ws.Rows("352:354").Insert()

# Common (constant) values shared by every row in this data block.
$mercadoId  = 9
$mercado    = "Vega Central Mapocho de Santiago"
$region     = "Metropolitana"
$codreg     = 13
$tipo       = "Fruta"
$productoId = 100107
$producto   = "Otros"
$categoriaId = 100107002
$categoria  = "Chirimoya"
$variedad   = "Cultivar IV Región"

$rows = @(
    @{ Row = 352; Fecha = 45212; Calidad = "Especial"; Volumen = 220; PMin = 27000; PMax = 27000; PProm = 27000; Unidad = "`$/bandeja 10 kilos"; Origen = "Provincia de Limarí"; PKg = 2700; KgUnidad = 10 },
    @{ Row = 353; Fecha = 45212; Calidad = "Primera";  Volumen = 290; PMin = 24000; PMax = 24000; PProm = 24000; Unidad = "`$/bandeja 10 kilos"; Origen = "Provincia de Limarí"; PKg = 2400; KgUnidad = 10 },
    @{ Row = 354; Fecha = 45212; Calidad = "Segunda";  Volumen = 280; PMin = 21000; PMax = 21000; PProm = 21000; Unidad = "`$/bandeja 10 kilos"; Origen = "Provincia de Limarí"; PKg = 2100; KgUnidad = 10 }
)

foreach ($r in $rows) {
    $i = $r.Row

    $ws.Cells.Item($i, 1).Value  = $mercadoId
    $ws.Cells.Item($i, 2).Value  = $mercado
    $ws.Cells.Item($i, 3).Value  = $region
    $ws.Cells.Item($i, 4).Value  = $r.Fecha
    $ws.Cells.Item($i, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($i, 5).Value  = $codreg
    $ws.Cells.Item($i, 6).Value  = $tipo
    $ws.Cells.Item($i, 7).Value  = $productoId
    $ws.Cells.Item($i, 8).Value  = $producto
    $ws.Cells.Item($i, 9).Value  = $categoriaId
    $ws.Cells.Item($i, 10).Value = $categoria
    $ws.Cells.Item($i, 11).Value = $variedad
    $ws.Cells.Item($i, 12).Value = $r.Calidad
    $ws.Cells.Item($i, 13).Value = $r.Volumen
    $ws.Cells.Item($i, 14).Value = $r.PMin
    $ws.Cells.Item($i, 15).Value = $r.PMax
    $ws.Cells.Item($i, 16).Value = $r.PProm
    $ws.Cells.Item($i, 17).Value = $r.Unidad
    $ws.Cells.Item($i, 18).Value = $r.Origen
    $ws.Cells.Item($i, 19).Value = $r.PKg
    $ws.Cells.Item($i, 20).Value = $r.KgUnidad
}
